$d = $word.ActiveDocument

# "Talk to User about customization options (e.g., color, any special
# requests, etc.)" -> "... (e.g., metal switch body, brackets, any special
# requests, etc.)" -- add the new metal foot pedal body / bracket option
# ahead of the existing "color" customization example.
$r = $d.Content
$r.Find.Execute("color", $true, $true, $false, $false, $false, $true, 1, $false, "metal switch body, brackets", 2)
